$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Before FS-DR")

# Row 75
$ws.Cells.Item(75, 2).Value = 0.9828
$ws.Cells.Item(75, 3).Value = 0.9825
$ws.Cells.Item(75, 4).Value = 0.9844000000000001
$ws.Cells.Item(75, 5).Value = 0.7294
$ws.Cells.Item(75, 6).Value = 1
$ws.Cells.Item(75, 7).Value = 1
$ws.Cells.Item(75, 8).Value = 1
$ws.Cells.Item(75, 9).Value = 1
$ws.Cells.Item(75, 10).Value = 0.9797
$ws.Cells.Item(75, 11).Value = 0.7302999999999999
$ws.Cells.Item(75, 12).Value = 'C=10'
$ws.Cells.Item(75, 13).Value = 'n_neighbors=3; weights=distance'
$ws.Cells.Item(75, 14).Value = 'max_depth=None; min_samples_split=2'
$ws.Cells.Item(75, 15).Value = 'max_depth=None; n_estimators=100'
$ws.Cells.Item(75, 16).Value = 'alpha=0.0001; hidden_layer_sizes=(50, 50)'

# Row 76
$ws.Cells.Item(76, 2).Value = 0.9797
$ws.Cells.Item(76, 3).Value = 0.9792
$ws.Cells.Item(76, 4).Value = 0.982
$ws.Cells.Item(76, 5).Value = 0.7211
$ws.Cells.Item(76, 6).Value = 1
$ws.Cells.Item(76, 7).Value = 1
$ws.Cells.Item(76, 8).Value = 1
$ws.Cells.Item(76, 9).Value = 1
$ws.Cells.Item(76, 10).Value = 0.9797
$ws.Cells.Item(76, 11).Value = 0.9734
$ws.Cells.Item(76, 12).Value = 'C=10'
$ws.Cells.Item(76, 13).Value = 'n_neighbors=3; weights=distance'
$ws.Cells.Item(76, 14).Value = 'max_depth=None; min_samples_split=2'
$ws.Cells.Item(76, 15).Value = 'max_depth=None; n_estimators=200'
$ws.Cells.Item(76, 16).Value = 'alpha=0.001; hidden_layer_sizes=(100,)'

# Row 77
$ws.Cells.Item(77, 2).Value = 0.982
$ws.Cells.Item(77, 3).Value = 0.9786
$ws.Cells.Item(77, 4).Value = 0.9859
$ws.Cells.Item(77, 5).Value = 0.731
$ws.Cells.Item(77, 6).Value = 1
$ws.Cells.Item(77, 7).Value = 1
$ws.Cells.Item(77, 8).Value = 0.9977
$ws.Cells.Item(77, 9).Value = 0.9977
$ws.Cells.Item(77, 10).Value = 0.9844000000000001
$ws.Cells.Item(77, 11).Value = 0.7318
$ws.Cells.Item(77, 12).Value = 'C=10'
$ws.Cells.Item(77, 13).Value = 'n_neighbors=3; weights=distance'
$ws.Cells.Item(77, 14).Value = 'max_depth=None; min_samples_split=2'
$ws.Cells.Item(77, 15).Value = 'max_depth=None; n_estimators=200'
$ws.Cells.Item(77, 16).Value = 'alpha=0.001; hidden_layer_sizes=(50, 50)'

# Row 78
$ws.Cells.Item(78, 2).Value = 0.9844000000000001
$ws.Cells.Item(78, 3).Value = 0.9842
$ws.Cells.Item(78, 4).Value = 0.9883
$ws.Cells.Item(78, 5).Value = 0.9784
$ws.Cells.Item(78, 6).Value = 1
$ws.Cells.Item(78, 7).Value = 1
$ws.Cells.Item(78, 8).Value = 0.9984
$ws.Cells.Item(78, 9).Value = 0.9985000000000001
$ws.Cells.Item(78, 10).Value = 0.9851
$ws.Cells.Item(78, 11).Value = 0.7327
$ws.Cells.Item(78, 12).Value = 'C=10'
$ws.Cells.Item(78, 13).Value = 'n_neighbors=3; weights=distance'
$ws.Cells.Item(78, 14).Value = 'max_depth=None; min_samples_split=2'
$ws.Cells.Item(78, 15).Value = 'max_depth=None; n_estimators=100'
$ws.Cells.Item(78, 16).Value = 'alpha=0.01; hidden_layer_sizes=(50, 50)'

# Row 79
$ws.Cells.Item(79, 2).Value = 0.9828
$ws.Cells.Item(79, 3).Value = 0.7333
$ws.Cells.Item(79, 4).Value = 0.9875
$ws.Cells.Item(79, 5).Value = 0.7383999999999999
$ws.Cells.Item(79, 6).Value = 1
$ws.Cells.Item(79, 7).Value = 1
$ws.Cells.Item(79, 8).Value = 1
$ws.Cells.Item(79, 9).Value = 1
$ws.Cells.Item(79, 10).Value = 0.9844000000000001
$ws.Cells.Item(79, 11).Value = 0.7314000000000001
$ws.Cells.Item(79, 12).Value = 'C=10'
$ws.Cells.Item(79, 13).Value = 'n_neighbors=5; weights=distance'
$ws.Cells.Item(79, 14).Value = 'max_depth=None; min_samples_split=2'
$ws.Cells.Item(79, 15).Value = 'max_depth=None; n_estimators=100'
$ws.Cells.Item(79, 16).Value = 'alpha=0.001; hidden_layer_sizes=(50, 50)'

# Row 80
$ws.Cells.Item(80, 2).Value = 0.9765
$ws.Cells.Item(80, 3).Value = 0.9694
$ws.Cells.Item(80, 4).Value = 0.9804
$ws.Cells.Item(80, 5).Value = 0.9707
$ws.Cells.Item(80, 6).Value = 1
$ws.Cells.Item(80, 7).Value = 1
$ws.Cells.Item(80, 8).Value = 0.9984
$ws.Cells.Item(80, 9).Value = 0.9985000000000001
$ws.Cells.Item(80, 10).Value = 0.9687
$ws.Cells.Item(80, 11).Value = 0.7
$ws.Cells.Item(80, 12).Value = 'C=10'
$ws.Cells.Item(80, 13).Value = 'n_neighbors=3; weights=distance'
$ws.Cells.Item(80, 14).Value = 'max_depth=None; min_samples_split=2'
$ws.Cells.Item(80, 15).Value = 'max_depth=None; n_estimators=100'
$ws.Cells.Item(80, 16).Value = 'alpha=0.0001; hidden_layer_sizes=(50, 50)'

# Row 81
$ws.Cells.Item(81, 2).Value = 0.9765
$ws.Cells.Item(81, 3).Value = 0.7229
$ws.Cells.Item(81, 4).Value = 0.9828
$ws.Cells.Item(81, 5).Value = 0.7272999999999999
$ws.Cells.Item(81, 6).Value = 1
$ws.Cells.Item(81, 7).Value = 1
$ws.Cells.Item(81, 8).Value = 0.9992
$ws.Cells.Item(81, 9).Value = 0.7463
$ws.Cells.Item(81, 10).Value = 0.9772999999999999
$ws.Cells.Item(81, 11).Value = 0.7238
$ws.Cells.Item(81, 12).Value = 'C=10'
$ws.Cells.Item(81, 13).Value = 'n_neighbors=3; weights=distance'
$ws.Cells.Item(81, 14).Value = 'max_depth=None; min_samples_split=2'
$ws.Cells.Item(81, 15).Value = 'max_depth=None; n_estimators=200'
$ws.Cells.Item(81, 16).Value = 'alpha=0.001; hidden_layer_sizes=(50, 50)'

# Row 82
$ws.Cells.Item(82, 2).Value = 0.9789
$ws.Cells.Item(82, 3).Value = 0.972
$ws.Cells.Item(82, 4).Value = 0.9789
$ws.Cells.Item(82, 5).Value = 0.7013
$ws.Cells.Item(82, 6).Value = 1
$ws.Cells.Item(82, 7).Value = 1
$ws.Cells.Item(82, 8).Value = 1
$ws.Cells.Item(82, 9).Value = 1
$ws.Cells.Item(82, 10).Value = 0.9703000000000001
$ws.Cells.Item(82, 11).Value = 0.7009
$ws.Cells.Item(82, 12).Value = 'C=10'
$ws.Cells.Item(82, 13).Value = 'n_neighbors=7; weights=distance'
$ws.Cells.Item(82, 14).Value = 'max_depth=None; min_samples_split=2'
$ws.Cells.Item(82, 15).Value = 'max_depth=None; n_estimators=100'
$ws.Cells.Item(82, 16).Value = 'alpha=0.001; hidden_layer_sizes=(50, 50)'

# Row 83
$ws.Cells.Item(83, 2).Value = 0.9828
$ws.Cells.Item(83, 3).Value = 0.7301
$ws.Cells.Item(83, 4).Value = 0.9867
$ws.Cells.Item(83, 5).Value = 0.7275
$ws.Cells.Item(83, 6).Value = 1
$ws.Cells.Item(83, 7).Value = 1
$ws.Cells.Item(83, 8).Value = 0.9992
$ws.Cells.Item(83, 9).Value = 0.9992
$ws.Cells.Item(83, 10).Value = 0.9796
$ws.Cells.Item(83, 11).Value = 0.7268
$ws.Cells.Item(83, 12).Value = 'C=10'
$ws.Cells.Item(83, 13).Value = 'n_neighbors=5; weights=distance'
$ws.Cells.Item(83, 14).Value = 'max_depth=None; min_samples_split=2'
$ws.Cells.Item(83, 15).Value = 'max_depth=None; n_estimators=200'
$ws.Cells.Item(83, 16).Value = 'alpha=0.0001; hidden_layer_sizes=(100,)'

# Row 84
$ws.Cells.Item(84, 2).Value = 0.9836
$ws.Cells.Item(84, 3).Value = 0.7341
$ws.Cells.Item(84, 4).Value = 0.9765
$ws.Cells.Item(84, 5).Value = 0.7208
$ws.Cells.Item(84, 6).Value = 1
$ws.Cells.Item(84, 7).Value = 1
$ws.Cells.Item(84, 8).Value = 1
$ws.Cells.Item(84, 9).Value = 1
$ws.Cells.Item(84, 10).Value = 0.9772999999999999
$ws.Cells.Item(84, 11).Value = 0.7276
$ws.Cells.Item(84, 12).Value = 'C=10'
$ws.Cells.Item(84, 13).Value = 'n_neighbors=3; weights=distance'
$ws.Cells.Item(84, 14).Value = 'max_depth=None; min_samples_split=2'
$ws.Cells.Item(84, 15).Value = 'max_depth=None; n_estimators=200'
$ws.Cells.Item(84, 16).Value = 'alpha=0.001; hidden_layer_sizes=(50, 50)'

# Row 87
$ws.Cells.Item(87, 2).Value = 0.9828
$ws.Cells.Item(87, 3).Value = 0.9825
$ws.Cells.Item(87, 4).Value = 0.9851
$ws.Cells.Item(87, 5).Value = 0.7302
$ws.Cells.Item(87, 6).Value = 1
$ws.Cells.Item(87, 7).Value = 1
$ws.Cells.Item(87, 8).Value = 1
$ws.Cells.Item(87, 9).Value = 1
$ws.Cells.Item(87, 10).Value = 0.9772999999999999
$ws.Cells.Item(87, 11).Value = 0.7275
$ws.Cells.Item(87, 12).Value = 'C=10'
$ws.Cells.Item(87, 13).Value = 'n_neighbors=3; weights=distance'
$ws.Cells.Item(87, 14).Value = 'max_depth=None; min_samples_split=2'
$ws.Cells.Item(87, 15).Value = 'max_depth=None; n_estimators=200'
$ws.Cells.Item(87, 16).Value = 'alpha=0.01; hidden_layer_sizes=(50,)'

# Row 88
$ws.Cells.Item(88, 2).Value = 0.9797
$ws.Cells.Item(88, 3).Value = 0.9792
$ws.Cells.Item(88, 4).Value = 0.982
$ws.Cells.Item(88, 5).Value = 0.7211
$ws.Cells.Item(88, 6).Value = 0.9992
$ws.Cells.Item(88, 7).Value = 0.9961
$ws.Cells.Item(88, 8).Value = 0.9992
$ws.Cells.Item(88, 9).Value = 0.9992
$ws.Cells.Item(88, 10).Value = 0.982
$ws.Cells.Item(88, 11).Value = 0.7296
$ws.Cells.Item(88, 12).Value = 'C=10'
$ws.Cells.Item(88, 13).Value = 'n_neighbors=3; weights=distance'
$ws.Cells.Item(88, 14).Value = 'max_depth=None; min_samples_split=2'
$ws.Cells.Item(88, 15).Value = 'max_depth=None; n_estimators=100'
$ws.Cells.Item(88, 16).Value = 'alpha=0.01; hidden_layer_sizes=(50, 50)'

# Row 89
$ws.Cells.Item(89, 2).Value = 0.982
$ws.Cells.Item(89, 3).Value = 0.9786
$ws.Cells.Item(89, 4).Value = 0.9859
$ws.Cells.Item(89, 5).Value = 0.731
$ws.Cells.Item(89, 6).Value = 0.9992
$ws.Cells.Item(89, 7).Value = 0.7499
$ws.Cells.Item(89, 8).Value = 0.9977
$ws.Cells.Item(89, 9).Value = 0.7484
$ws.Cells.Item(89, 10).Value = 0.9875
$ws.Cells.Item(89, 11).Value = 0.735
$ws.Cells.Item(89, 12).Value = 'C=10'
$ws.Cells.Item(89, 13).Value = 'n_neighbors=3; weights=distance'
$ws.Cells.Item(89, 14).Value = 'max_depth=None; min_samples_split=2'
$ws.Cells.Item(89, 15).Value = 'max_depth=None; n_estimators=100'
$ws.Cells.Item(89, 16).Value = 'alpha=0.0001; hidden_layer_sizes=(50, 50)'

# Row 90
$ws.Cells.Item(90, 2).Value = 0.9836
$ws.Cells.Item(90, 3).Value = 0.9834000000000001
$ws.Cells.Item(90, 4).Value = 0.9883
$ws.Cells.Item(90, 5).Value = 0.9784
$ws.Cells.Item(90, 6).Value = 1
$ws.Cells.Item(90, 7).Value = 1
$ws.Cells.Item(90, 8).Value = 0.9984
$ws.Cells.Item(90, 9).Value = 0.9985000000000001
$ws.Cells.Item(90, 10).Value = 0.9851
$ws.Cells.Item(90, 11).Value = 0.7327
$ws.Cells.Item(90, 12).Value = 'C=10'
$ws.Cells.Item(90, 13).Value = 'n_neighbors=3; weights=distance'
$ws.Cells.Item(90, 14).Value = 'max_depth=None; min_samples_split=5'
$ws.Cells.Item(90, 15).Value = 'max_depth=None; n_estimators=100'
$ws.Cells.Item(90, 16).Value = 'alpha=0.001; hidden_layer_sizes=(50, 50)'

# Row 91
$ws.Cells.Item(91, 2).Value = 0.9828
$ws.Cells.Item(91, 3).Value = 0.7333
$ws.Cells.Item(91, 4).Value = 0.9875
$ws.Cells.Item(91, 5).Value = 0.7383
$ws.Cells.Item(91, 6).Value = 1
$ws.Cells.Item(91, 7).Value = 1
$ws.Cells.Item(91, 8).Value = 1
$ws.Cells.Item(91, 9).Value = 1
$ws.Cells.Item(91, 10).Value = 0.9844000000000001
$ws.Cells.Item(91, 11).Value = 0.7314000000000001
$ws.Cells.Item(91, 12).Value = 'C=10'
$ws.Cells.Item(91, 13).Value = 'n_neighbors=5; weights=distance'
$ws.Cells.Item(91, 14).Value = 'max_depth=None; min_samples_split=10'
$ws.Cells.Item(91, 15).Value = 'max_depth=20; n_estimators=100'
$ws.Cells.Item(91, 16).Value = 'alpha=0.001; hidden_layer_sizes=(50, 50)'

# Row 92
$ws.Cells.Item(92, 2).Value = 0.9765
$ws.Cells.Item(92, 3).Value = 0.9694
$ws.Cells.Item(92, 4).Value = 0.9797
$ws.Cells.Item(92, 5).Value = 0.97
$ws.Cells.Item(92, 6).Value = 1
$ws.Cells.Item(92, 7).Value = 1
$ws.Cells.Item(92, 8).Value = 0.9984
$ws.Cells.Item(92, 9).Value = 0.9985000000000001
$ws.Cells.Item(92, 10).Value = 0.9844000000000001
$ws.Cells.Item(92, 11).Value = 0.8918
$ws.Cells.Item(92, 12).Value = 'C=10'
$ws.Cells.Item(92, 13).Value = 'n_neighbors=3; weights=distance'
$ws.Cells.Item(92, 14).Value = 'max_depth=10; min_samples_split=2'
$ws.Cells.Item(92, 15).Value = 'max_depth=None; n_estimators=100'
$ws.Cells.Item(92, 16).Value = 'alpha=0.001; hidden_layer_sizes=(50, 50)'

# Row 93
$ws.Cells.Item(93, 2).Value = 0.9765
$ws.Cells.Item(93, 3).Value = 0.7229
$ws.Cells.Item(93, 4).Value = 0.9836
$ws.Cells.Item(93, 5).Value = 0.7281
$ws.Cells.Item(93, 6).Value = 0.9992
$ws.Cells.Item(93, 7).Value = 0.7463
$ws.Cells.Item(93, 8).Value = 0.9984
$ws.Cells.Item(93, 9).Value = 0.7455000000000001
$ws.Cells.Item(93, 10).Value = 0.975
$ws.Cells.Item(93, 11).Value = 0.7183
$ws.Cells.Item(93, 12).Value = 'C=10'
$ws.Cells.Item(93, 13).Value = 'n_neighbors=3; weights=distance'
$ws.Cells.Item(93, 14).Value = 'max_depth=None; min_samples_split=2'
$ws.Cells.Item(93, 15).Value = 'max_depth=None; n_estimators=100'
$ws.Cells.Item(93, 16).Value = 'alpha=0.0001; hidden_layer_sizes=(100,)'

# Row 94
$ws.Cells.Item(94, 2).Value = 0.9789
$ws.Cells.Item(94, 3).Value = 0.972
$ws.Cells.Item(94, 4).Value = 0.9789
$ws.Cells.Item(94, 5).Value = 0.6986
$ws.Cells.Item(94, 6).Value = 1
$ws.Cells.Item(94, 7).Value = 1
$ws.Cells.Item(94, 8).Value = 1
$ws.Cells.Item(94, 9).Value = 1
$ws.Cells.Item(94, 10).Value = 0.9656
$ws.Cells.Item(94, 11).Value = 0.6833
$ws.Cells.Item(94, 12).Value = 'C=10'
$ws.Cells.Item(94, 13).Value = 'n_neighbors=7; weights=distance'
$ws.Cells.Item(94, 14).Value = 'max_depth=None; min_samples_split=2'
$ws.Cells.Item(94, 15).Value = 'max_depth=None; n_estimators=200'
$ws.Cells.Item(94, 16).Value = 'alpha=0.01; hidden_layer_sizes=(50, 50)'

# Row 95
$ws.Cells.Item(95, 2).Value = 0.9828
$ws.Cells.Item(95, 3).Value = 0.7301
$ws.Cells.Item(95, 4).Value = 0.9867
$ws.Cells.Item(95, 5).Value = 0.7275
$ws.Cells.Item(95, 6).Value = 1
$ws.Cells.Item(95, 7).Value = 1
$ws.Cells.Item(95, 8).Value = 0.9977
$ws.Cells.Item(95, 9).Value = 0.9977
$ws.Cells.Item(95, 10).Value = 0.9843
$ws.Cells.Item(95, 11).Value = 0.7285
$ws.Cells.Item(95, 12).Value = 'C=10'
$ws.Cells.Item(95, 13).Value = 'n_neighbors=5; weights=distance'
$ws.Cells.Item(95, 14).Value = 'max_depth=None; min_samples_split=5'
$ws.Cells.Item(95, 15).Value = 'max_depth=None; n_estimators=100'
$ws.Cells.Item(95, 16).Value = 'alpha=0.0001; hidden_layer_sizes=(50, 50)'

# Row 96
$ws.Cells.Item(96, 2).Value = 0.982
$ws.Cells.Item(96, 3).Value = 0.7324000000000001
$ws.Cells.Item(96, 4).Value = 0.9765
$ws.Cells.Item(96, 5).Value = 0.724
$ws.Cells.Item(96, 6).Value = 1
$ws.Cells.Item(96, 7).Value = 1
$ws.Cells.Item(96, 8).Value = 0.9984
$ws.Cells.Item(96, 9).Value = 0.9985000000000001
$ws.Cells.Item(96, 10).Value = 0.9812
$ws.Cells.Item(96, 11).Value = 0.7319
$ws.Cells.Item(96, 12).Value = 'C=10'
$ws.Cells.Item(96, 13).Value = 'n_neighbors=3; weights=distance'
$ws.Cells.Item(96, 14).Value = 'max_depth=None; min_samples_split=2'
$ws.Cells.Item(96, 15).Value = 'max_depth=None; n_estimators=100'
$ws.Cells.Item(96, 16).Value = 'alpha=0.01; hidden_layer_sizes=(50, 50)'

